$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first three product rows (CAPIXY HAIR CREAM / SHAMPOO / TONIC SPRAY),
# shifting the remaining rows up.
$ws.Rows("7:9").Delete()

# Renumber the serial ("م") column for the remaining 7 rows.
for ($i = 0; $i -lt 7; $i++) {
    $row = 7 + $i
    $ws.Cells.Item($row, 1).Value = $i + 1
}

# Update the displayed grand total to match the remaining rows.
$ws.Range("N14").Value = 249.595
